# Add CVSCC results, and fix some spelling errors in COM results.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add CVSCC result -------------------------------------------------------
# Row 40 (XS-B / BLAKE DOBLINGER) gets a CVSCC ("Z") final_time result.
$ws.Range("A40").Value = "Z"
$ws.Range("E40").Value = " 0:56.329"

# --- Fix spelling errors ---------------------------------------------------
# Row 47: "jJOHN MORAVEC" -> "JOHN MORAVEC"
$ws.Range("C47").Value = " JOHN MORAVEC"

# Row 28: "KIM JOHM CRUMB" -> "KIM JOHN CRUMB"
$ws.Range("C28").Value = " KIM JOHN CRUMB"

# Row 9: "DENNIS OMEARA" -> "DENNIS O'MEARA"
$ws.Range("C9").Value = " DENNIS O'MEARA"
